# Swap the contents of columns B:AC between two pairs of rows,
# leaving column A (the sequential row id) untouched.
# Uses .Value2 cell-by-cell (per-cell) since whole-range .Value
# assignment is not reliable in this COM-interop runtime.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-Rows($Worksheet, $Row1, $Row2, $Columns) {
    $vals1 = @{}
    $vals2 = @{}

    foreach ($col in $Columns) {
        $vals1[$col] = $Worksheet.Range("$col$Row1").Value2
        $vals2[$col] = $Worksheet.Range("$col$Row2").Value2
    }

    foreach ($col in $Columns) {
        $Worksheet.Range("$col$Row1").Value2 = $vals2[$col]
        $Worksheet.Range("$col$Row2").Value2 = $vals1[$col]
    }
}

# Rows 135 and 136: swap columns B through AC
Swap-Rows $ws 135 136 $cols

# Rows 152 and 153: swap columns B through AC
Swap-Rows $ws 152 153 $cols
